# Add data for 2022-10-01
# This updates the "carjacking by neighborhood by month" report to extend
# the current (partial) month window from "through September 22" to
# "through September 23", and bumps the year-over-year September counts
# for the neighborhoods that had a carjacking recorded on 09-23 of each
# year (plus the running current-month total).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet / title text changes.
$ws.Name = "Through 2022-09-23"
$ws.Range("B1").Value = "September 2022 (through September 23)"

# Cell value updates (row => column => new value).
$ws.Range("B3").Value = 7
$ws.Range("AC3").Value = 4
$ws.Range("AU3").Value = 5

$ws.Range("BM9").Value = 1

$ws.Range("B10").Value = 6

$ws.Range("AU12").Value = 2

$ws.Range("AL14").Value = 3

$ws.Range("AU26").Value = 2

$ws.Range("AL27").Value = 1
$ws.Range("AU27").Value = 2

$ws.Range("T33").Value = 4

$ws.Range("T34").Value = 1

$ws.Range("K43").Value = 2

$ws.Range("AL48").Value = 1

$ws.Range("AC50").Value = 2

$ws.Range("AU56").Value = 1

$ws.Range("K63").Value = 2

$ws.Range("AU64").Value = 2

$ws.Range("BM73").Value = 1

$ws.Range("T77").Value = 2

$ws.Range("B91").Value = 2

$ws.Range("B96").Value = 4

$ws.Range("AL98").Value = 3
